$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on B:E for rows 2-51 so that numeric-looking
# strings (e.g. "174.24", "0.998") are stored as text, matching the source data.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "67.386.34"
$ws.Cells.Item(2, 5).Value = "  -0.13%  "

$ws.Cells.Item(3, 4).Value = "3.283.58"
$ws.Cells.Item(3, 5).Value = "  -2.73%  "

$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  +0.26%  "

$ws.Cells.Item(5, 4).Value = "580.88"
$ws.Cells.Item(5, 5).Value = "  -1.72%  "

$ws.Cells.Item(6, 4).Value = "174.24"
$ws.Cells.Item(6, 5).Value = "  -7.66%  "

$ws.Cells.Item(7, 4).Value = "0.998"
$ws.Cells.Item(7, 5).Value = "  -0.06%  "

$ws.Cells.Item(8, 5).Value = "  -3.51%  "

$ws.Cells.Item(9, 4).Value = "3.276.16"
$ws.Cells.Item(9, 5).Value = "  -2.89%  "

$ws.Cells.Item(10, 4).Value = "0.173"
$ws.Cells.Item(10, 5).Value = "  -5.92%  "

$ws.Cells.Item(11, 4).Value = "0.569"
$ws.Cells.Item(11, 5).Value = "  -3.37%  "

$ws.Cells.Item(12, 4).Value = "44.98"
$ws.Cells.Item(12, 5).Value = "  -5.61%  "

$ws.Cells.Item(13, 4).Value = "0.0000267"
$ws.Cells.Item(13, 5).Value = "  -2.88%  "

$ws.Cells.Item(14, 4).Value = "669.28"
$ws.Cells.Item(14, 5).Value = "  +4.72%  "

$ws.Cells.Item(15, 4).Value = "3.808.26"
$ws.Cells.Item(15, 5).Value = "  -2.73%  "

$ws.Cells.Item(16, 4).Value = "8.28"
$ws.Cells.Item(16, 5).Value = "  -4.21%  "

$ws.Cells.Item(17, 4).Value = "67.391.31"
$ws.Cells.Item(17, 5).Value = "  -0.01%  "

$ws.Cells.Item(18, 5).Value = "  -0.44%  "

$ws.Cells.Item(19, 4).Value = "3.284.65"
$ws.Cells.Item(19, 5).Value = "  -2.43%  "

$ws.Cells.Item(20, 4).Value = "17.30"
$ws.Cells.Item(20, 5).Value = "  -4.33%  "

$ws.Cells.Item(21, 4).Value = "10.77"
$ws.Cells.Item(21, 5).Value = "  -4.04%  "

$ws.Cells.Item(22, 4).Value = "0.880"
$ws.Cells.Item(22, 5).Value = "  -3.64%  "

$ws.Cells.Item(23, 2).Value = "Toncoin"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(23, 4).Value = "5.40"
$ws.Cells.Item(23, 5).Value = "  +5.71%  "

$ws.Cells.Item(24, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(24, 4).Value = "16.98"
$ws.Cells.Item(24, 5).Value = "  -6.11%  "

$ws.Cells.Item(25, 4).Value = "97.70"
$ws.Cells.Item(25, 5).Value = "  -2.87%  "

$ws.Cells.Item(26, 4).Value = "3.85"
$ws.Cells.Item(26, 5).Value = "  -4.27%  "

$ws.Cells.Item(27, 4).Value = "2.64"
$ws.Cells.Item(27, 5).Value = "  -7.26%  "

$ws.Cells.Item(28, 4).Value = "9.14"
$ws.Cells.Item(28, 5).Value = "  -6.57%  "

$ws.Cells.Item(29, 4).Value = "32.70"
$ws.Cells.Item(29, 5).Value = "  +0.06%  "

$ws.Cells.Item(30, 4).Value = "8.30"
$ws.Cells.Item(30, 5).Value = "  -5.01%  "

$ws.Cells.Item(31, 4).Value = "6.89"
$ws.Cells.Item(31, 5).Value = "  +0.56%  "

$ws.Cells.Item(32, 4).Value = "569.47"
$ws.Cells.Item(32, 5).Value = "  -7.28%  "

$ws.Cells.Item(33, 4).Value = "10.87"
$ws.Cells.Item(33, 5).Value = "  -2.99%  "

$ws.Cells.Item(34, 4).Value = "3.740.06"
$ws.Cells.Item(34, 5).Value = "  -4.75%  "

$ws.Cells.Item(35, 4).Value = "0.102"
$ws.Cells.Item(35, 5).Value = "  -4.26%  "

$ws.Cells.Item(36, 4).Value = "0.998"
$ws.Cells.Item(36, 5).Value = "  -0.34%  "

$ws.Cells.Item(37, 5).Value = "  -12.64%  "

$ws.Cells.Item(38, 4).Value = "55.66"
$ws.Cells.Item(38, 5).Value = "  -0.65%  "

$ws.Cells.Item(39, 5).Value = "  -1.69%  "

$ws.Cells.Item(40, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(40, 4).Value = "32.26"
$ws.Cells.Item(40, 5).Value = "  -4.80%  "

$ws.Cells.Item(41, 2).Value = "Fetch.AI"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(41, 4).Value = "2.61"
$ws.Cells.Item(41, 5).Value = "  -8.22%  "

$ws.Cells.Item(42, 4).Value = "3.03"
$ws.Cells.Item(42, 5).Value = "  -7.11%  "

$ws.Cells.Item(43, 4).Value = "0.0₃0658"
$ws.Cells.Item(43, 5).Value = "  -7.13%  "

$ws.Cells.Item(44, 2).Value = "ApeXProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(44, 4).Value = "3.22"
$ws.Cells.Item(44, 5).Value = "  -5.40%  "

$ws.Cells.Item(45, 2).Value = "TheGraph"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(45, 4).Value = "0.325"
$ws.Cells.Item(45, 5).Value = "  -6.16%  "

$ws.Cells.Item(46, 4).Value = "0.0400"
$ws.Cells.Item(46, 5).Value = "  -5.59%  "

$ws.Cells.Item(47, 4).Value = "2.58"
$ws.Cells.Item(47, 5).Value = "  -0.57%  "

$ws.Cells.Item(48, 5).Value = "  +0.18%  "

$ws.Cells.Item(49, 5).Value = "  -2.80%  "

$ws.Cells.Item(50, 4).Value = "1.35"
$ws.Cells.Item(50, 5).Value = "  -0.76%  "

$ws.Cells.Item(51, 4).Value = "2.75"
$ws.Cells.Item(51, 5).Value = "  -4.22%  "

# Restore the original (default/"Normal") style so no stray number-format
# styling is left behind on the cells.
$ws.Range("B2:E51").Style = "Normal"
